$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches source data which is
# always stored as inline/shared strings, even for number-looking values
# like prices "1.00" or "0.468" that Excel would otherwise auto-convert
# to numbers and strip formatting from).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '59.215.26'
$ws.Range('E2').Value = '  -7.57%  '

Set-TextValue $ws.Range('D3') '3.305.94'
$ws.Range('E3').Value = '  -4.79%  '

Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  -0.01%  '

Set-TextValue $ws.Range('D5') '559.88'
$ws.Range('E5').Value = '  -4.17%  '

Set-TextValue $ws.Range('D6') '127.14'
$ws.Range('E6').Value = '  -2.73%  '

$ws.Range('E7').Value = '  +0.04%  '

Set-TextValue $ws.Range('D8') '3.306.92'
$ws.Range('E8').Value = '  -4.71%  '

Set-TextValue $ws.Range('D9') '0.468'
$ws.Range('E9').Value = '  -2.91%  '

Set-TextValue $ws.Range('D10') '7.34'
$ws.Range('E10').Value = '  -4.28%  '

$ws.Range('E11').Value = '  -5.88%  '

$ws.Range('E12').Value = '  -3.73%  '

Set-TextValue $ws.Range('D13') '3.876.06'
$ws.Range('E13').Value = '  -4.78%  '

$ws.Range('E14').Value = '  -0.16%  '

Set-TextValue $ws.Range('D15') '3.318.83'
$ws.Range('E15').Value = '  -4.32%  '

$ws.Range('E16').Value = '  -6.49%  '

Set-TextValue $ws.Range('D17') '23.92'
$ws.Range('E17').Value = '  -4.08%  '

Set-TextValue $ws.Range('D18') '59.536.27'

Set-TextValue $ws.Range('D19') '5.61'
$ws.Range('E19').Value = '  -0.94%  '

Set-TextValue $ws.Range('D20') '13.20'
$ws.Range('E20').Value = '  -1.15%  '

Set-TextValue $ws.Range('D21') '8.85'
$ws.Range('E21').Value = '  -10.96%  '

Set-TextValue $ws.Range('D22') '348.99'
$ws.Range('E22').Value = '  -9.54%  '

Set-TextValue $ws.Range('D23') '0.551'
$ws.Range('E23').Value = '  -2.43%  '

$ws.Range('E24').Value = '  +0.21%  '

Set-TextValue $ws.Range('D25') '3.439.82'
$ws.Range('E25').Value = '  -4.83%  '

Set-TextValue $ws.Range('D26') '68.34'
$ws.Range('E26').Value = '  -8.21%  '

$ws.Range('E27').Value = '  -3.07%  '

Set-TextValue $ws.Range('D28') '1.00'
$ws.Range('E28').Value = '  +0.04%  '

Set-TextValue $ws.Range('D29') '7.26'
$ws.Range('E29').Value = '  +2.76%  '

Set-TextValue $ws.Range('D30') '1.43'
$ws.Range('E30').Value = '  +0.23%  '

Set-TextValue $ws.Range('D31') '7.75'
$ws.Range('E31').Value = '  -2.29%  '

$ws.Range('E32').Value = '  -3.94%  '

$ws.Range('E33').Value = '  -6.08%  '

$ws.Range('E34').Value = '  +0.02%  '

Set-TextValue $ws.Range('D35') '3.339.88'
$ws.Range('E35').Value = '  -4.64%  '

Set-TextValue $ws.Range('D36') '22.67'
$ws.Range('E36').Value = '  -1.27%  '

Set-TextValue $ws.Range('D37') '5.20'
$ws.Range('E37').Value = '  +0.32%  '

Set-TextValue $ws.Range('D38') '6.74'
$ws.Range('E38').Value = '  -0.09%  '

Set-TextValue $ws.Range('D39') '1.46'
$ws.Range('E39').Value = '  -1.89%  '

Set-TextValue $ws.Range('D40') '157.33'
$ws.Range('E40').Value = '  -3.51%  '

Set-TextValue $ws.Range('D41') '0.0742'
$ws.Range('E41').Value = '  -3.89%  '

Set-TextValue $ws.Range('D42') '1.00'
$ws.Range('E42').Value = '  -0.05%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D43') '0.741'
$ws.Range('E43').Value = '  -6.78%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D44') '40.29'
$ws.Range('E44').Value = '  -2.95%  '

$ws.Range('E45').Value = '  -1.62%  '

$ws.Range('E46').Value = '  +3.56%  '

Set-TextValue $ws.Range('D47') '22.68'
$ws.Range('E47').Value = '  -3.32%  '

Set-TextValue $ws.Range('D48') '1.52'
$ws.Range('E48').Value = '  -5.73%  '

Set-TextValue $ws.Range('D49') '6.68'
$ws.Range('E49').Value = '  -0.18%  '

Set-TextValue $ws.Range('D50') '21.61'
$ws.Range('E50').Value = '  +6.02%  '

$ws.Range('E51').Value = '  +11.49%  '
